# CLI Execution.xlsx — remove duplicated CLI command lines that were
# accidentally left in from copy/paste of the vpls configuration blocks.
#
# Sheet "172.31.72.93": a duplicated "no shutdown" line (the second one of a
# repeated "no shutdown"/"exit" pair) is removed from each of the four vpls
# blocks (originally rows 22, 39, 59, 76).
#
# Sheet "172.31.72.94": both a stray "tempna" line and a duplicated
# "no shutdown" line are removed from each of the four vpls blocks
# (originally rows 9, 17, 24, 35, 42, 53, 60, 68).
#
# Rows are removed highest-row-number first so earlier row numbers keep
# referring to the same original rows while we work through the list.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("172.31.72.93")
$sheet1Rows = @(76, 59, 39, 22)
foreach ($r in $sheet1Rows) {
    $ws1.Rows.Item($r).Delete()
}

$ws2 = $wb.Worksheets.Item("172.31.72.94")
$sheet2Rows = @(68, 60, 53, 42, 35, 24, 17, 9)
foreach ($r in $sheet2Rows) {
    $ws2.Rows.Item($r).Delete()
}
